# Update "tasas-transfi.xlsx" with the latest automated rate values.

$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the summary text in A1 -------------------------
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.74 = 10267.24 pesos`n✅ 10267.24 pesos = 2.73 = 940.22 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update the rate figures --------------------------------
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 364.99
$wsTasas.Range("O10").Value = 3747.44

$wsTasas.Range("N12").Value = 3765
$wsTasas.Range("O12").Value = 344.78
